# Update Handback status report timestamps (Generate Report for Handback)

$wb = $excel.ActiveWorkbook

# Overview sheet - "Latest HO Xliff Generate Date" for first file row
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 01:10:19"

# zh-cn sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 01:10:11"
$wsZhCn.Range("K2").Value = "2016-09-01 01:10:38"

# de-de sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 01:10:19"
$wsDeDe.Range("K2").Value = "2016-09-01 01:10:46"
